$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns: D = porcentaje_utilidades, E = porcentaje_contingencia
$ws.Range("D1").Value = "porcentaje_utilidades"
$ws.Range("E1").Value = "porcentaje_contingencia"
$ws.Range("E2").Value = 13

# Header formatting for the two new header cells (bold, centered, thin left/right border)
foreach ($addr in @("D1", "E1")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.Item(7).Weight = 2
    $cell.Borders.Item(10).Weight = 2
}

# Column widths to fit the new headers
$ws.Columns.Item(3).ColumnWidth = 21.453125
$ws.Columns.Item(4).ColumnWidth = 22.08984375
$ws.Columns.Item(5).ColumnWidth = 24.6328125

# Restore cursor/selection like the saved file shows
$null = $ws.Range("E4").Select()
